$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 10 into the new row 11 so the new label cell (A11)
# picks up the same bold/border/centered style used by the other period labels.
$ws.Range("A10:G10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2 (Q0) ---
$ws.Range("B2").Value = 0.003902179103014683
$ws.Range("C2").Value = 0.4110727462040676
$ws.Range("D2").Value = 0.3428415174671495
$ws.Range("E2").Value = 0.5855267009002659
$ws.Range("F2").Value = 0.5913398483737633
$ws.Range("G2").Value = 51

# --- Row 3 (Q1) ---
$ws.Range("B3").Value = 0.1114323317728629
$ws.Range("C3").Value = 0.3819837740696647
$ws.Range("D3").Value = 0.3070771422459426
$ws.Range("E3").Value = 0.5541454161553108
$ws.Range("F3").Value = 0.548336982715814
$ws.Range("G3").Value = 50

# --- Row 4 (Q2) ---
$ws.Range("B4").Value = 0.0253973765991602
$ws.Range("C4").Value = 0.4132096390943262
$ws.Range("D4").Value = 0.334517719201454
$ws.Range("E4").Value = 0.5783750679286357
$ws.Range("F4").Value = 0.5838050818178278
$ws.Range("G4").Value = 49

# --- Row 5 (Q3) ---
$ws.Range("B5").Value = 0.1183269686296411
$ws.Range("C5").Value = 0.3749468204795861
$ws.Range("D5").Value = 0.3172515219839373
$ws.Range("E5").Value = 0.563250851738315
$ws.Range("F5").Value = 0.5565091045686217
$ws.Range("G5").Value = 48

# --- Row 6 (Q4) ---
$ws.Range("B6").Value = 0.03052796929569665
$ws.Range("C6").Value = 0.3923170207885465
$ws.Range("D6").Value = 0.3392183943349846
$ws.Range("E6").Value = 0.5824245825297767
$ws.Range("F6").Value = 0.5879119750550927
$ws.Range("G6").Value = 47

# --- Row 7 (Q5) ---
$ws.Range("B7").Value = 0.1229440676987344
$ws.Range("C7").Value = 0.3639192002280563
$ws.Range("D7").Value = 0.3218065997724393
$ws.Range("E7").Value = 0.5672800012096666
$ws.Range("F7").Value = 0.5599167076062078
$ws.Range("G7").Value = 46

# --- Row 8 (Q6) ---
$ws.Range("B8").Value = 0.02678724024085173
$ws.Range("C8").Value = 0.3918890428263422
$ws.Range("D8").Value = 0.3443878750560778
$ws.Range("E8").Value = 0.5868456995293378
$ws.Range("F8").Value = 0.5928583371096039
$ws.Range("G8").Value = 45

# --- Row 9 (Q7) ---
$ws.Range("B9").Value = 0.06383516716290022
$ws.Range("C9").Value = 0.3005467402051563
$ws.Range("D9").Value = 0.2218349939727257
$ws.Range("E9").Value = 0.4709936241317134
$ws.Range("F9").Value = 0.4720426389359249
$ws.Range("G9").Value = 44

# --- Row 10 (Q8) --- (note: F10 previously had no value at all, now it does)
$ws.Range("B10").Value = 0.07455416172597187
$ws.Range("C10").Value = 0.3804291969162187
$ws.Range("D10").Value = 0.3199665900932944
$ws.Range("E10").Value = 0.5656558937139208
$ws.Range("F10").Value = 0.5673571874781979
$ws.Range("G10").Value = 43

# --- Row 11 (Q9) --- NEW ROW
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.04160261366545656
$ws.Range("C11").Value = 0.3112244733113839
$ws.Range("D11").Value = 0.2248588685484449
$ws.Range("E11").Value = 0.4741928600774635
$ws.Range("F11").Value = 0.4780902003259299
$ws.Range("G11").Value = 42
